$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update state/zip values for the two records (new service area parsed from XML -> ON postal codes)
$ws.Range("E2").Value = "L4L4Y8"
$ws.Range("D2").Value = "ON"
$ws.Range("D3").Value = "ON"
$ws.Range("E3").Value = "L7E4G4"

# Autosize the data columns (A:E) to fit their content (matches Excel's
# AutoFit-computed "best fit" widths for the Name/address/city/state/zip columns)
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(2).ColumnWidth = 17.833333333333332
$ws.Columns.Item(3).ColumnWidth = 9.666666666666666
$ws.Columns.Item(4).ColumnWidth = 4.666666666666667
$ws.Columns.Item(5).ColumnWidth = 6.166666666666667

# Move/update the active selection
$ws.Range("F3").Select()
